# Applies the "fix report acdc project & menu hrd" update:
#  - Updates report period titles from Dec/22 to Jan/23 on both sheets
#  - Appends transaction rows (3-5) to the "Transaction Maker ACDC" sheet
#  - Appends project summary rows (5-7) to the "Worksheet" sheet
#
# Helper to write a value into a cell while forcing it to stay as TEXT,
# even when it looks like a date (e.g. "2023-01-18"), without leaving any
# stray NumberFormat/quotePrefix styling behind on the cell (plain
# "$cell.Value = ..." lets Excel auto-detect such strings as dates and
# permanently stamps a new cell style on them). We build it as a text
# formula first (so Excel can't reinterpret it), then collapse the
# formula down to its plain cached value via copy / paste-values.
function Set-TextCell {
    param($cell, $text)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Transaction Maker ACDC"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Update report period in the title (Dec/22 -> Jan/23)
$ws1.Cells.Item(1,1).Value = "Export Laporan Transaction Maker ACDC Periode 01/Jan/23 - 31/Jan/23"

# Row 3
$ws1.Cells.Item(3,1).Value = "Maintain"
$ws1.Cells.Item(3,2).Value = 1234565

# Row 4
$ws1.Cells.Item(4,1).Value = "Bea Cukai"
$ws1.Cells.Item(4,2).Value = 123213
$ws1.Cells.Item(4,3).Value = "Transfer"
$ws1.Cells.Item(4,4).Value = "Jono"
$ws1.Cells.Item(4,5).Value = 2000000
$ws1.Cells.Item(4,6).Value = "Beli permen"
Set-TextCell $ws1.Cells.Item(4,7) "2023-01-18"

# Row 5
$ws1.Cells.Item(5,1).Value = "Renewal Redhat"
$ws1.Cells.Item(5,2).Value = 80000123
$ws1.Cells.Item(5,3).Value = "Payment"
$ws1.Cells.Item(5,4).Value = "Jon1"
$ws1.Cells.Item(5,5).Value = 10000000
$ws1.Cells.Item(5,6).Value = "beli permen"
Set-TextCell $ws1.Cells.Item(5,7) "2023-01-13"

# ---------------------------------------------------------------------
# Sheet 2: "Worksheet"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Update report period in the title (Dec/22 -> Jan/23)
$ws2.Cells.Item(1,1).Value = "Export Laporan Data Project ACDC Periode 01/Jan/23 - 31/Jan/23"

# Row 5
$ws2.Cells.Item(5,1).Value = 1234565
$ws2.Cells.Item(5,2).Value = "Maintain"
$ws2.Cells.Item(5,3).Value = "BC"
$ws2.Cells.Item(5,4).Value = "AP"
$ws2.Cells.Item(5,5).Value = 90000000
$ws2.Cells.Item(5,6).Value = 1000000
Set-TextCell $ws2.Cells.Item(5,7) "1.000.000"
$ws2.Cells.Item(5,8).Value = 92000000
$ws2.Cells.Item(5,9).Value = 4600000
$ws2.Cells.Item(5,10).Value = 5
$ws2.Cells.Item(5,11).Value = 20000000
$ws2.Cells.Item(5,12).Value = 67400000
$ws2.Cells.Item(5,13).Value = "2023-01-19 20:12:08"

# Row 6
$ws2.Cells.Item(6,1).Value = 123213
$ws2.Cells.Item(6,2).Value = "Bea Cukai"
$ws2.Cells.Item(6,3).Value = "D1"
$ws2.Cells.Item(6,4).Value = "AP"
$ws2.Cells.Item(6,5).Value = 800000000
$ws2.Cells.Item(6,6).Value = 250000000
Set-TextCell $ws2.Cells.Item(6,7) "1.000.000"
$ws2.Cells.Item(6,8).Value = 1051000000
$ws2.Cells.Item(6,9).Value = 84080000
$ws2.Cells.Item(6,10).Value = 8
$ws2.Cells.Item(6,11).Value = 6000000
$ws2.Cells.Item(6,12).Value = 960920000
$ws2.Cells.Item(6,13).Value = "2023-01-22 01:49:16"

# Row 7
$ws2.Cells.Item(7,1).Value = 80000123
$ws2.Cells.Item(7,2).Value = "Renewal Redhat"
$ws2.Cells.Item(7,3).Value = "D1"
$ws2.Cells.Item(7,4).Value = "AP"
$ws2.Cells.Item(7,5).Value = 800000000
$ws2.Cells.Item(7,6).Value = 250000000
Set-TextCell $ws2.Cells.Item(7,7) "1.000.000"
$ws2.Cells.Item(7,8).Value = 1051000000
$ws2.Cells.Item(7,9).Value = 84080000
$ws2.Cells.Item(7,10).Value = 8
$ws2.Cells.Item(7,11).Value = 6000000
$ws2.Cells.Item(7,12).Value = 960920000
$ws2.Cells.Item(7,13).Value = "2023-01-22 01:52:35"
